# Find.Execute positional signature used throughout this script:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#           MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
# Replace must be 2 (wdReplaceAll) to actually commit every match inside the
# searched range/story, or 1 (wdReplaceOne) to commit just the first match.

$d = $word.ActiveDocument

# 1) Title (Heading1) + the bolded repeat near the bottom of the document.
#    Both occurrences are identical, so one document-wide wdReplaceAll
#    Execute call updates both runs while preserving each run's own
#    formatting (bold stays bold, heading stays a heading).
$d.Content.Find.Execute(
    "Play Blown Away Slot for Free - Review & Bonuses",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Play Blown Away Free - Exciting Animal-Themed Slot Game",
    2) | Out-Null

# 2) "What we like" bullets - each phrase is unique in the document.
$d.Content.Find.Execute(
    "Engaging animal-inspired theme",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Engaging gameplay with animal-inspired theme",
    2) | Out-Null

$d.Content.Find.Execute(
    "Various bonus features",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Various bonuses to boost chances of winning prizes",
    2) | Out-Null

$d.Content.Find.Execute(
    "Visually appealing graphics",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Visually appealing and immersive graphics",
    2) | Out-Null

# "25 always active paylines" appears several times in the body copy, but
# only the standalone bullet-list paragraph should gain the extra phrase.
# Walk the paragraphs to find the exact (whole-paragraph) bullet text, then
# run Find.Execute scoped to just that one paragraph's Range so the other
# sentences that merely contain this phrase are left untouched.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "25 always active paylines`r") {
        $scoped = $d.Range($para.Range.Start, $para.Range.End)
        $scoped.Find.Execute(
            "25 always active paylines",
            $true, $false, $false, $false, $false,
            $true, 1, $false,
            "25 always active paylines for increased winning potential",
            2) | Out-Null
        break
    }
}

# 3) "What we don't like" bullet: split into two bullets. Putting "^p" in
#    the ReplaceWith text inserts a real paragraph mark (this special
#    replacement code is recognised regardless of MatchWildcards, same as
#    in real Word), so the single bullet becomes two ListBullet paragraphs
#    -- both inheriting the original paragraph's formatting -- in one
#    Execute call.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Medium volatility with lower RTP compared to industry average`r") {
        $scoped = $d.Range($para.Range.Start, $para.Range.End)
        $scoped.Find.Execute(
            "Medium volatility with lower RTP compared to industry average",
            $true, $false, $false, $false, $false,
            $true, 1, $false,
            "Medium volatility with better paying options available^pSlightly lower RTP compared to industry average",
            2) | Out-Null
        break
    }
}

# 4) Meta description (italic paragraph at the end of the document).
$d.Content.Find.Execute(
    "Explore the unique animal-inspired theme of Blown Away slot and enjoy various bonus features. Play for free and learn about the game's RTP and volatility.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Blown Away, an engaging slot game with animal characters. Play for free and win prizes!",
    2) | Out-Null
